$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 106, shifting the existing rows 106-108 down to 108-110
$ws.Rows.Item(106).Insert()
$ws.Rows.Item(106).Insert()

# New row 106: Black Amber, "Primera"
$ws.Cells.Item(106,1).Value = 4
$ws.Cells.Item(106,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(106,3).Value = "Los Lagos"
$ws.Cells.Item(106,4).Value = 44568
$ws.Cells.Item(106,5).Value = 10
$ws.Cells.Item(106,6).Value = "Fruta"
$ws.Cells.Item(106,7).Value = 100103
$ws.Cells.Item(106,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(106,9).Value = 100103002
$ws.Cells.Item(106,10).Value = "Ciruela"
$ws.Cells.Item(106,11).Value = "Black Amber"
$ws.Cells.Item(106,12).Value = "Primera"
$ws.Cells.Item(106,13).Value = 600
$ws.Cells.Item(106,14).Value = 18000
$ws.Cells.Item(106,15).Value = 18500
$ws.Cells.Item(106,16).Value = 18250
$ws.Cells.Item(106,17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(106,18).Value = "Región de O'Higgins"
$ws.Cells.Item(106,19).Value = 1217
$ws.Cells.Item(106,20).Value = 15

# New row 107: Black Amber, "Segunda"
$ws.Cells.Item(107,1).Value = 4
$ws.Cells.Item(107,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(107,3).Value = "Los Lagos"
$ws.Cells.Item(107,4).Value = 44568
$ws.Cells.Item(107,5).Value = 10
$ws.Cells.Item(107,6).Value = "Fruta"
$ws.Cells.Item(107,7).Value = 100103
$ws.Cells.Item(107,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(107,9).Value = 100103002
$ws.Cells.Item(107,10).Value = "Ciruela"
$ws.Cells.Item(107,11).Value = "Black Amber"
$ws.Cells.Item(107,12).Value = "Segunda"
$ws.Cells.Item(107,13).Value = 200
$ws.Cells.Item(107,14).Value = 16000
$ws.Cells.Item(107,15).Value = 16000
$ws.Cells.Item(107,16).Value = 16000
$ws.Cells.Item(107,17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(107,18).Value = "Región de O'Higgins"
$ws.Cells.Item(107,19).Value = 1067
$ws.Cells.Item(107,20).Value = 15
